$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.400.24'
$ws.Range('E2').Value = '  -6.60%  '

$ws.Range('D3').Value = '2.936.26'
$ws.Range('E3').Value = '  -8.61%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '531.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -10.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -15.92%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = '2.904.19'
$ws.Range('E8').Value = '  -9.31%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.447'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -17.43%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.141'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -17.87%  '

$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.80'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.92%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.423'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -14.06%  '

$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '31.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -19.45%  '

$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000198'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -18.18%  '

$ws.Range('D15').Value = '3.412.59'
$ws.Range('E15').Value = '  -8.67%  '

$ws.Range('D16').Value = '62.195.96'
$ws.Range('E16').Value = '  -7.01%  '

$ws.Range('E17').Value = '  -4.66%  '

$ws.Range('D18').Value = '2.925.54'
$ws.Range('E18').Value = '  -9.10%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '465.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.88%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -14.93%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -16.51%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.621'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -17.78%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -20.54%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -13.66%  '

$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -14.85%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -18.10%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.81'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -15.64%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -17.29%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '24.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -16.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.83%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.994'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.95%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -15.22%  '

$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.72'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.53%  '

$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '460.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -15.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.31'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -17.26%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.56'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -19.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0378'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -10.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0741'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.09%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.108'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.92%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -18.09%  '

$ws.Range('D42').Value = '2.586.73'
$ws.Range('E42').Value = '  -10.77%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.997'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.29%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -18.07%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.215'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -18.38%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '109.76'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.73%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.14%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.76'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -18.27%  '

$ws.Range('D49').Value = '0.0₃0454'
$ws.Range('E49').Value = '  -20.96%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.91'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -20.50%  '

$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.15'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.57%  '
